# Update "想去人数" (F column) and one "最低票价" (G15) value on the
# "展览" and "全部类型" sheets to reflect the newer snapshot of counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 1040
    $ws.Range("F4").Value = 515
    $ws.Range("F5").Value = 13699
    $ws.Range("F6").Value = 42
    $ws.Range("F7").Value = 35
    $ws.Range("F8").Value = 1756
    $ws.Range("F9").Value = 161
    $ws.Range("F15").Value = 13695
    $ws.Range("G15").Value = 60
    $ws.Range("F18").Value = 9036
    $ws.Range("F20").Value = 8155
    $ws.Range("F28").Value = 1025

    if ($sheetName -eq "展览") {
        $ws.Range("F32").Value = 4
        $ws.Range("F37").Value = 2967
    } else {
        $ws.Range("F34").Value = 4
        $ws.Range("F39").Value = 2968
    }
}
